$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.494.72"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").Value = "1.846.96"
$ws.Range("E3").Value = "  -1.36%  "
$ws.Range("D4").Formula = "=""0.9992"""
$ws.Range("D4").Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Formula = "=""261.66"""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -7.20%  "
$ws.Range("D6").Formula = "=""0.9999"""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").Formula = "=""0.5129"""
$ws.Range("D7").Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("D8").Formula = "=""0.3216"""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -9.16%  "
$ws.Range("D9").Formula = "=""0.06770"""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -1.73%  "
$ws.Range("D10").Formula = "=""19.00"""
$ws.Range("D10").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  -5.81%  "
$ws.Range("D11").Formula = "=""0.7689"""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -5.90%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Formula = "=""0.07696"""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.860.77"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").Formula = "=""88.96"""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("D15").Formula = "=""5.033"""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("D16").Formula = "=""0.9986"""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Formula = "=""14.11"""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -1.94%  "
$ws.Range("D18").Formula = "=""0.9989"""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Formula = "=""0.000007910"""
$ws.Range("D19").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("D20").Value = "26.519.34"
$ws.Range("E20").Value = "  -0.60%  "
$ws.Range("D21").Value = "2.111.88"
$ws.Range("E21").Value = "  +0.30%  "
$ws.Range("E22").Value = "  -5.04%  "
$ws.Range("D23").Formula = "=""9.554"""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -5.65%  "
$ws.Range("D24").Formula = "=""5.964"""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -4.40%  "
$ws.Range("D25").Formula = "=""2.348"""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").Formula = "=""144.87"""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Formula = "=""1.661"""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  -2.06%  "
$ws.Range("D29").Formula = "=""111.18"""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.38%  "
$ws.Range("D30").Formula = "=""4.210"""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("D31").Formula = "=""4.176"""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -3.82%  "
$ws.Range("D32").Formula = "=""0.08730"""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("D33").Formula = "=""0.04832"""
$ws.Range("D33").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("D34").Formula = "=""1.140"""
$ws.Range("D34").Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  -3.11%  "
$ws.Range("D35").Formula = "=""2.846"""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").Formula = "=""0.6906"""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -7.44%  "
$ws.Range("D37").Formula = "=""3.099"""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -5.47%  "
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("D39").Formula = "=""2.211"""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -8.27%  "
$ws.Range("D40").Formula = "=""0.4918"""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -6.16%  "
$ws.Range("D41").Formula = "=""113.78"""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").Formula = "=""0.9047"""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -6.78%  "
$ws.Range("D43").Formula = "=""6.158"""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -2.21%  "
$ws.Range("D44").Formula = "=""0.9991"""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +0.34%  "
$ws.Range("D45").Formula = "=""7.807"""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -4.16%  "
$ws.Range("D46").Formula = "=""0.4259"""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -7.22%  "
$ws.Range("D47").Formula = "=""0.1271"""
$ws.Range("D47").Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -6.84%  "
$ws.Range("D48").Formula = "=""9.137"""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  -3.90%  "
$ws.Range("D49").Formula = "=""0.05896"""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -0.39%  "
$ws.Range("D50").Formula = "=""34.99"""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Formula = "=""1.425"""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -5.87%  "
$excel.CutCopyMode = $false
